# T1503_Contact_ContactDetails_AddEditDeleteActivity.xlsx
# Mid - 16th Jan 2025
#
# Update the "Contact" sheet sample row (row 3): the contact name moves
# from the ad-hoc "Test Houlihan" to the real Houlihan Lokey employee
# type, and the related-company value is shortened to "HL".
# Also move the saved selection on that sheet to A9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

$ws.Range("A3").Value = "Houlihan Employee"
$ws.Range("B3").Value = "HL"

$ws.Activate()
$ws.Range("A9").Select()
